$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the new timesheet entry for row 4 (matching style of rows 2 and 3)
$ws.Range("B4").Value = 8
$ws.Range("C4").Value = "Improved FileOpener readability and efficiency."

# Update the selected cell to C4 (as seen in the diff's sheetView selection)
$ws.Range("C4").Select()

# Recalculate so the SUM formula in B38 reflects the new total
$excel.Calculate()
